$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (raw xlsx "width" = ColumnWidth + 5/6, so subtract
# the 0.8333... padding offset to land exactly on the target stored width)
$ws.Columns.Item(1).ColumnWidth = 88.16666666666667
$ws.Columns.Item(3).ColumnWidth = 30.166666666666668
$ws.Columns.Item(4).ColumnWidth = 32.166666666666664
$ws.Columns.Item(5).ColumnWidth = 39.166666666666664
$ws.Columns.Item(6).ColumnWidth = 41.166666666666664

# Update header labels (row 1)
$ws.Range("C1").Value = "label_storageColorOptions_for"
$ws.Range("D1").Value = "label_storageColorOptions_for_1"
$ws.Range("E1").Value = "label_storageColorOptions_internalText"
$ws.Range("F1").Value = "label_storageColorOptions_internalText_1"

# Update data path value in row 2
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/fillShippingAddressDetailsAndContinueToPayment-test-data"
